# "Clean up code and fix output"
# Adds a new "Yearly demand" worksheet (as the last tab) that holds the
# hourly yearly-demand profile table: a header row of hour indices 0-23
# (columns B:Y), a left-hand column of day-type indices 0-2 (A2:A4), and
# the corresponding demand values in the B2:Y4 body.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last tab ("Connected Households")
# so it lands at the end of the workbook, and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Yearly demand"

# Header row: hours 0-23 across B1:Y1.
$headerRow = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23)
for ($c = 0; $c -lt $headerRow.Length; $c++) {
    $ws.Cells.Item(1, $c + 2).Value = $headerRow[$c]
}

# Left column: day-type index 0-2 across A2:A4.
$colA = @(0, 1, 2)
for ($r = 0; $r -lt $colA.Length; $r++) {
    $ws.Cells.Item($r + 2, 1).Value = $colA[$r]
}

# Data body B2:Y4 - hourly demand values for each of the 3 day types.
$bodyData = @(
    @(-32.5, -19.5, -13, -13, -13, 142.5, 291.5, 327, 388.5, 502, 596, 670.5, 745, 651, 576.5, 502, 320.5, 139, 32, -117, -97.5, -78, -52, -39),
    @(-32.5, -19.5, -13, 0, 0, -19.5, 0, 324, 486, 648, 729, 751.5, 583, 567, 333.5, 340, 243, 57.99999999999999, -130, 0, 0, -78, 0, -39),
    @(-32.5, -19.5, 0, 0, 0, -19.5, 0, 0, 81, 324, 567, 589.5, 648, 567, 324, 162, 81, 0, -130, 0, 0, 0, 0, -39)
)
for ($r = 0; $r -lt $bodyData.Length; $r++) {
    $rowVals = $bodyData[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}

# Style the header row and the day-type column like the rest of the
# workbook's table headers: bold, centered/top-aligned, thin box border.
$headerRange = $ws.Range("B1:Y1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$dayTypeRange = $ws.Range("A2:A4")
$dayTypeRange.Font.Bold = $true
$dayTypeRange.HorizontalAlignment = -4108
$dayTypeRange.VerticalAlignment = -4160
$dayTypeRange.Borders.LineStyle = 1

# Leave selection on A1 like the other sheets in this workbook.
$ws.Range("A1").Select()
